$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the two "failed login" test e-mail addresses (3rd and 4th rows)
$ws.Range("A3").Value = "test@basarisiz.com"
$ws.Range("A4").Value = "basarisiz@basarisiz.com"

# Resize column A (best-fit) so the longer e-mail addresses fit
$ws.Columns.Item(1).ColumnWidth = 21.8

# Move the active selection to B5
$ws.Range("B5").Select() | Out-Null
